$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Meta")
$ws = $wb.Worksheets.Item("Instructions")
$ws.Range("C1").Value = "Template updated 1/17/23"
$ws.Range("C1").Font.Color = 255
$ws.Activate()
$ws.Range("B15").Select()
$meta.Activate()
